$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round row 5 (B5:AH5) values to 2 decimal places ---------------------
$row5Values = @(
    20.34, 15.59, 0.58, 43.19, 36.88, 15.6, 57.56,
    24.14, 10.93, 16.71, 17.84, 18.83, 5.03, 15.8,
    22.28, 12.71, 0.48, 0.73, 232.84, 43.77, 14.79,
    30.09, 15.65, 2.07, 28.3, 12.56, 11.41, 13.39,
    18.88, 0, 51.63, 8.45, 18.06
)
for ($i = 0; $i -lt $row5Values.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5Values[$i]
}

# --- Delete row 6 entirely -------------------------------------------------
$ws.Rows(6).Delete()

# --- Narrow a subset of data columns from 8 (or 9) chars to 7 (or 8) -------
$columnWidths = @{
    2  = 7;  3  = 7;  7  = 7;  9  = 7;  10 = 7; 11 = 7; 13 = 7; 15 = 7
    16 = 7;  17 = 7;  20 = 8;  22 = 7;  23 = 7; 24 = 7; 26 = 7; 27 = 7
    28 = 7;  29 = 7;  34 = 7
}
foreach ($col in $columnWidths.Keys) {
    $ws.Cells.Item(1, $col).EntireColumn.ColumnWidth = $columnWidths[$col] - 5 / 6
}
